# Insert two new weekly price rows for "Apio" (Terminal La Palmera de La Serena)
# right before the existing row 564, shifting all subsequent rows down by two.
# The new rows capture the most recent week's "Primera" / "Segunda" quality data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 564 (pushes old 564.. down to 566..)
$ws.Rows("564:565").Insert()

# New row 564: Primera
$ws.Range("A564").Value = 8
$ws.Range("B564").Value = "Terminal La Palmera de La Serena"
$ws.Range("C564").Value = "Coquimbo"
$ws.Range("D564").Value = 45077
$ws.Range("E564").Value = 4
$ws.Range("F564").Value = 100112017
$ws.Range("G564").Value = "Apio"
$ws.Range("H564").Value = "Americana (o)"
$ws.Range("I564").Value = "Primera"
$ws.Range("J564").Value = 1600
$ws.Range("K564").Value = 8000
$ws.Range("L564").Value = 9000
$ws.Range("M564").Value = 8500
$ws.Range("N564").Value = "`$/docena de matas"
$ws.Range("O564").Value = "Provincia del Elquí"
$ws.Range("P564").Value = 1417
$ws.Range("Q564").Value = 6
$ws.Range("R564").Value = "Hortaliza"

# New row 565: Segunda
$ws.Range("A565").Value = 8
$ws.Range("B565").Value = "Terminal La Palmera de La Serena"
$ws.Range("C565").Value = "Coquimbo"
$ws.Range("D565").Value = 45077
$ws.Range("E565").Value = 4
$ws.Range("F565").Value = 100112017
$ws.Range("G565").Value = "Apio"
$ws.Range("H565").Value = "Americana (o)"
$ws.Range("I565").Value = "Segunda"
$ws.Range("J565").Value = 800
$ws.Range("K565").Value = 6000
$ws.Range("L565").Value = 7000
$ws.Range("M565").Value = 6500
$ws.Range("N565").Value = "`$/docena de matas"
$ws.Range("O565").Value = "Provincia del Elquí"
$ws.Range("P565").Value = 1083
$ws.Range("Q565").Value = 6
$ws.Range("R565").Value = "Hortaliza"
